$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 2615.3333  # H17: 2778.4 -> 2615.3333
$ws.Cells.Item(17, 10).Value = 2615.3333  # J17: 2778.4 -> 2615.3333
$ws.Cells.Item(17, 12).Value = 7845.999899999999  # L17: 8335.200000000001 -> 7845.999899999999
$ws.Cells.Item(17, 14).Value = -8181.999899999999  # N17: -8671.200000000001 -> -8181.999899999999
$ws.Cells.Item(19, 8).Value = 480.3889  # H19: 469.94736 -> 480.3889
$ws.Cells.Item(19, 10).Value = 658.3  # J19: 624.0909 -> 658.3
$ws.Cells.Item(19, 12).Value = 658.3  # L19: 624.0909 -> 658.3
$ws.Cells.Item(19, 14).Value = -1008.3  # N19: -974.0909 -> -1008.3
$ws.Cells.Item(40, 8).Value = 4959.2856  # H40: 4891.4546 -> 4959.2856
$ws.Cells.Item(40, 9).Value = 4281.154  # I40: 4223 -> 4281.154
$ws.Cells.Item(40, 11).Value = 4281.154  # K40: 4223 -> 4281.154
$ws.Cells.Item(40, 13).Value = -4106.154  # M40: -4048 -> -4106.154
$ws.Cells.Item(51, 8).Value = 10813.2  # H51: 10813.3 -> 10813.2
$ws.Cells.Item(51, 9).Value = 10675  # I51: 12333.667 -> 10675
$ws.Cells.Item(51, 10).Value = 10905.333  # J51: 10161.714 -> 10905.333
$ws.Cells.Item(51, 11).Value = 10675  # K51: 12333.667 -> 10675
$ws.Cells.Item(51, 12).Value = 10905.333  # L51: 10161.714 -> 10905.333
$ws.Cells.Item(51, 13).Value = -10191  # M51: -11849.667 -> -10191
$ws.Cells.Item(51, 14).Value = -11873.333  # N51: -11129.714 -> -11873.333
$ws.Cells.Item(109, 8).Value = 105000  # H109: 150000 -> 105000
$ws.Cells.Item(109, 10).Value = 105000  # J109: 150000 -> 105000
$ws.Cells.Item(109, 12).Value = 105000  # L109: 150000 -> 105000
$ws.Cells.Item(109, 14).Value = -107774  # N109: -152774 -> -107774
$ws.Cells.Item(137, 8).Value = 2178.262  # H137: 2007.7843 -> 2178.262
$ws.Cells.Item(137, 9).Value = 2162.2222  # I137: 1845.5555 -> 2162.2222
$ws.Cells.Item(137, 11).Value = 6486.6666  # K137: 5536.666499999999 -> 6486.6666
$ws.Cells.Item(137, 13).Value = -3936.6666  # M137: -2986.666499999999 -> -3936.6666
$ws.Cells.Item(138, 8).Value = 3347.59  # H138: 253958.8 -> 3347.59
$ws.Cells.Item(138, 9).Value = 2711.75  # I138: 4531.7827 -> 2711.75
$ws.Cells.Item(138, 10).Value = 3506.55  # J138: 328462.97 -> 3506.55
$ws.Cells.Item(138, 11).Value = 8135.25  # K138: 13595.3481 -> 8135.25
$ws.Cells.Item(138, 12).Value = 10519.65  # L138: 985388.9099999999 -> 10519.65
$ws.Cells.Item(138, 13).Value = -2995.25  # M138: -8455.348099999999 -> -2995.25
$ws.Cells.Item(138, 14).Value = -20799.65  # N138: -995668.9099999999 -> -20799.65
$ws.Cells.Item(141, 8).Value = 692.7143  # H141: 774.8333 -> 692.7143
$ws.Cells.Item(141, 10).Value = 862.5  # J141: 1083.3334 -> 862.5
$ws.Cells.Item(141, 12).Value = 2587.5  # L141: 3250.0002 -> 2587.5
$ws.Cells.Item(141, 14).Value = -12947.5  # N141: -13610.0002 -> -12947.5

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 1885.7468  # H32: 1680.81 -> 1885.7468
$ws.Cells.Item(32, 9).Value = 1885.7468  # I32: 1694.7677 -> 1885.7468
$ws.Cells.Item(32, 10).Value = 0  # J32: 299 -> 0
$ws.Cells.Item(32, 11).Value = 1885.7468  # K32: 1694.7677 -> 1885.7468
$ws.Cells.Item(32, 12).Value = 0  # L32: 299 -> 0
$ws.Cells.Item(32, 13).ClearContents()  # M32: delete (was -1407.7677)
$ws.Cells.Item(32, 14).Value = -1598.7468  # N32: -873 -> -1598.7468
$ws.Cells.Item(61, 8).Value = 5306.88  # H61: 5313.923 -> 5306.88
$ws.Cells.Item(61, 9).Value = 4333.647  # I61: 4397.8887 -> 4333.647
$ws.Cells.Item(61, 11).Value = 4333.647  # K61: 4397.8887 -> 4333.647
$ws.Cells.Item(61, 13).Value = -4121.647  # M61: -4185.8887 -> -4121.647
$ws.Cells.Item(122, 8).Value = 3053.1167  # H122: 3054.8667 -> 3053.1167
$ws.Cells.Item(122, 9).Value = 2915  # I122: 2917.0588 -> 2915
$ws.Cells.Item(122, 11).Value = 8745  # K122: 8751.1764 -> 8745
$ws.Cells.Item(122, 13).Value = -6295  # M122: -6301.1764 -> -6295
$ws.Cells.Item(136, 8).Value = 5306.88  # H136: 5313.923 -> 5306.88
$ws.Cells.Item(136, 9).Value = 4333.647  # I136: 4397.8887 -> 4333.647
$ws.Cells.Item(136, 11).Value = 13000.941  # K136: 13193.6661 -> 13000.941
$ws.Cells.Item(136, 13).Value = -10450.941  # M136: -10643.6661 -> -10450.941

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value = 3979.2222  # H99: 4066.2354 -> 3979.2222
$ws.Cells.Item(99, 10).Value = 5496.6665  # J99: 5871.25 -> 5496.6665
$ws.Cells.Item(99, 12).Value = 5496.6665  # L99: 5871.25 -> 5496.6665
$ws.Cells.Item(99, 14).Value = -8492.666499999999  # N99: -8867.25 -> -8492.666499999999
$ws.Cells.Item(134, 8).Value = 4060  # H134: 3729.25 -> 4060
$ws.Cells.Item(134, 9).Value = 3920.2  # I134: 3260.4 -> 3920.2
$ws.Cells.Item(134, 10).Value = 4199.8  # J134: 4064.1428 -> 4199.8
$ws.Cells.Item(134, 11).Value = 11760.6  # K134: 9781.200000000001 -> 11760.6
$ws.Cells.Item(134, 12).Value = 12599.4  # L134: 12192.4284 -> 12599.4
$ws.Cells.Item(134, 13).Value = -9225.599999999999  # M134: -7246.200000000001 -> -9225.599999999999
$ws.Cells.Item(134, 14).Value = -17669.4  # N134: -17262.4284 -> -17669.4

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2765.2  # H31: 3424.3572 -> 2765.2
$ws.Cells.Item(31, 9).Value = 2468.2856  # I31: 2919.139 -> 2468.2856
$ws.Cells.Item(31, 10).Value = 3759.2173  # J31: 4823.423 -> 3759.2173
$ws.Cells.Item(31, 11).Value = 2468.2856  # K31: 2919.139 -> 2468.2856
$ws.Cells.Item(31, 12).Value = 3759.2173  # L31: 4823.423 -> 3759.2173
$ws.Cells.Item(31, 13).Value = -2173.2856  # M31: -2624.139 -> -2173.2856
$ws.Cells.Item(31, 14).Value = -4349.2173  # N31: -5413.423 -> -4349.2173
$ws.Cells.Item(34, 8).Value = 2765.2  # H34: 3424.3572 -> 2765.2
$ws.Cells.Item(34, 9).Value = 2468.2856  # I34: 2919.139 -> 2468.2856
$ws.Cells.Item(34, 10).Value = 3759.2173  # J34: 4823.423 -> 3759.2173
$ws.Cells.Item(34, 11).Value = 2468.2856  # K34: 2919.139 -> 2468.2856
$ws.Cells.Item(34, 12).Value = 3759.2173  # L34: 4823.423 -> 3759.2173
$ws.Cells.Item(34, 13).Value = -2266.2856  # M34: -2717.139 -> -2266.2856
$ws.Cells.Item(34, 14).Value = -4163.2173  # N34: -5227.423 -> -4163.2173
$ws.Cells.Item(58, 8).Value = 3273.9565  # H58: 3350.0454 -> 3273.9565
$ws.Cells.Item(58, 9).Value = 2395.8  # I58: 2484.2222 -> 2395.8
$ws.Cells.Item(58, 11).Value = 2395.8  # K58: 2484.2222 -> 2395.8
$ws.Cells.Item(58, 13).Value = -2192.8  # M58: -2281.2222 -> -2192.8
$ws.Cells.Item(107, 8).Value = 845.7  # H107: 883.8421 -> 845.7
$ws.Cells.Item(107, 9).Value = 882.125  # I107: 932.86664 -> 882.125
$ws.Cells.Item(107, 11).Value = 882.125  # K107: 932.86664 -> 882.125
$ws.Cells.Item(107, 13).Value = 1037.875  # M107: 987.13336 -> 1037.875
$ws.Cells.Item(136, 8).Value = 3273.9565  # H136: 3350.0454 -> 3273.9565
$ws.Cells.Item(136, 9).Value = 2395.8  # I136: 2484.2222 -> 2395.8
$ws.Cells.Item(136, 11).Value = 7187.400000000001  # K136: 7452.6666 -> 7187.400000000001
$ws.Cells.Item(136, 13).Value = -4637.400000000001  # M136: -4902.6666 -> -4637.400000000001
$ws.Cells.Item(141, 8).Value = 571620  # H141: 341759 -> 571620
$ws.Cells.Item(141, 10).Value = 571620  # J141: 341759 -> 571620
$ws.Cells.Item(141, 12).Value = 571620  # L141: 341759 -> 571620
$ws.Cells.Item(141, 14).Value = -581980  # N141: -352119 -> -581980

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 2443.3  # H113: 1892.5385 -> 2443.3
$ws.Cells.Item(113, 9).Value = 1750  # I113: 506.66666 -> 1750
$ws.Cells.Item(113, 10).Value = 2520.3333  # J113: 2308.3 -> 2520.3333
$ws.Cells.Item(113, 11).Value = 5250  # K113: 1519.99998 -> 5250
$ws.Cells.Item(113, 12).Value = 7560.999899999999  # L113: 6924.900000000001 -> 7560.999899999999
$ws.Cells.Item(113, 13).Value = -3080  # M113: 650.0000199999999 -> -3080
$ws.Cells.Item(113, 14).Value = -11900.9999  # N113: -11264.9 -> -11900.9999
$ws.Cells.Item(122, 8).Value = 1383.9231  # H122: 1415.9166 -> 1383.9231
$ws.Cells.Item(122, 10).Value = 1364.6364  # J122: 1401.1 -> 1364.6364
$ws.Cells.Item(122, 12).Value = 12281.7276  # L122: 12609.9 -> 12281.7276
$ws.Cells.Item(122, 14).Value = -17181.7276  # N122: -17509.9 -> -17181.7276

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2318.476  # H102: 2151.4 -> 2318.476
$ws.Cells.Item(102, 9).Value = 2383.75  # I102: 2037.625 -> 2383.75
$ws.Cells.Item(102, 10).Value = 2231.4443  # J102: 2353.6667 -> 2231.4443
$ws.Cells.Item(102, 11).Value = 2383.75  # K102: 2037.625 -> 2383.75
$ws.Cells.Item(102, 12).Value = 2231.4443  # L102: 2353.6667 -> 2231.4443
$ws.Cells.Item(102, 13).Value = -761.75  # M102: -415.625 -> -761.75
$ws.Cells.Item(102, 14).Value = -5475.4443  # N102: -5597.6667 -> -5475.4443
$ws.Cells.Item(126, 8).Value = 7215.067  # H126: 7215 -> 7215.067
$ws.Cells.Item(126, 9).Value = 3201.2  # I126: 3201.1 -> 3201.2
$ws.Cells.Item(126, 11).Value = 9603.599999999999  # K126: 9603.299999999999 -> 9603.599999999999
$ws.Cells.Item(126, 13).Value = -7133.599999999999  # M126: -7133.299999999999 -> -7133.599999999999

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(4, 8).Value = 19998  # H4: 19997.334 -> 19998
$ws.Cells.Item(4, 10).Value = 19998  # J4: 19997.334 -> 19998
$ws.Cells.Item(4, 12).Value = 19998  # L4: 19997.334 -> 19998
$ws.Cells.Item(4, 14).Value = -20224  # N4: -20223.334 -> -20224
$ws.Cells.Item(22, 8).Value = 76532730  # H22: 107145150 -> 76532730
$ws.Cells.Item(22, 9).Value = 7144467  # I22: 14287343 -> 7144467
$ws.Cells.Item(22, 10).Value = 250003380  # J22: 200002960 -> 250003380
$ws.Cells.Item(22, 11).Value = 7144467  # K22: 14287343 -> 7144467
$ws.Cells.Item(22, 12).Value = 250003380  # L22: 200002960 -> 250003380
$ws.Cells.Item(22, 13).Value = -7144172  # M22: -14287048 -> -7144172
$ws.Cells.Item(22, 14).Value = -250003970  # N22: -200003550 -> -250003970
$ws.Cells.Item(27, 8).Value = 76532730  # H27: 107145150 -> 76532730
$ws.Cells.Item(27, 9).Value = 7144467  # I27: 14287343 -> 7144467
$ws.Cells.Item(27, 10).Value = 250003380  # J27: 200002960 -> 250003380
$ws.Cells.Item(27, 11).Value = 7144467  # K27: 14287343 -> 7144467
$ws.Cells.Item(27, 12).Value = 250003380  # L27: 200002960 -> 250003380
$ws.Cells.Item(27, 13).Value = -7144360  # M27: -14287236 -> -7144360
$ws.Cells.Item(27, 14).Value = -250003594  # N27: -200003174 -> -250003594
$ws.Cells.Item(28, 8).Value = 19998  # H28: 19997.334 -> 19998
$ws.Cells.Item(28, 10).Value = 19998  # J28: 19997.334 -> 19998
$ws.Cells.Item(28, 12).Value = 19998  # L28: 19997.334 -> 19998
$ws.Cells.Item(28, 14).Value = -20462  # N28: -20461.334 -> -20462
$ws.Cells.Item(37, 8).Value = 19998  # H37: 19997.334 -> 19998
$ws.Cells.Item(37, 10).Value = 19998  # J37: 19997.334 -> 19998
$ws.Cells.Item(37, 12).Value = 19998  # L37: 19997.334 -> 19998
$ws.Cells.Item(37, 14).Value = -20212  # N37: -20211.334 -> -20212
$ws.Cells.Item(100, 8).Value = 5796.1055  # H100: 5869.278 -> 5796.1055
$ws.Cells.Item(100, 9).Value = 4845.4375  # I100: 4869.8667 -> 4845.4375
$ws.Cells.Item(100, 11).Value = 4845.4375  # K100: 4869.8667 -> 4845.4375
$ws.Cells.Item(100, 13).Value = -4304.4375  # M100: -4328.8667 -> -4304.4375
$ws.Cells.Item(112, 8).Value = 64387  # H112: 0 -> 64387
$ws.Cells.Item(112, 10).Value = 64387  # J112: 0 -> 64387
$ws.Cells.Item(112, 12).Value = 64387  # L112: 0 -> 64387
$ws.Cells.Item(112, 14).Value = -67341  # N112: None -> -67341
$ws.Cells.Item(132, 8).Value = 2991.2917  # H132: 3827.2 -> 2991.2917
$ws.Cells.Item(132, 9).Value = 2389.7  # I132: 2605.5625 -> 2389.7
$ws.Cells.Item(132, 10).Value = 5999.25  # J132: 5999 -> 5999.25
$ws.Cells.Item(132, 11).Value = 7169.099999999999  # K132: 7816.6875 -> 7169.099999999999
$ws.Cells.Item(132, 12).Value = 17997.75  # L132: 17997 -> 17997.75
$ws.Cells.Item(132, 13).Value = -4639.099999999999  # M132: -5286.6875 -> -4639.099999999999
$ws.Cells.Item(132, 14).Value = -23057.75  # N132: -23057 -> -23057.75
$ws.Cells.Item(136, 8).Value = 4653.645  # H136: 6388.7 -> 4653.645
$ws.Cells.Item(136, 9).Value = 4139.227  # I136: 4403.55 -> 4139.227
$ws.Cells.Item(136, 10).Value = 5911.1113  # J136: 10359 -> 5911.1113
$ws.Cells.Item(136, 11).Value = 12417.681  # K136: 13210.65 -> 12417.681
$ws.Cells.Item(136, 12).Value = 17733.3339  # L136: 31077 -> 17733.3339
$ws.Cells.Item(136, 13).Value = -9867.681  # M136: -10660.65 -> -9867.681
$ws.Cells.Item(136, 14).Value = -22833.3339  # N136: -36177 -> -22833.3339

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(26, 8).Value = 15000  # H26: 0 -> 15000
$ws.Cells.Item(26, 9).Value = 15000  # I26: 0 -> 15000
$ws.Cells.Item(26, 11).Value = 15000  # K26: 0 -> 15000
$ws.Cells.Item(26, 13).Value = -14707  # M26: None -> -14707
$ws.Cells.Item(81, 8).Value = 6723  # H81: 4820.737 -> 6723
$ws.Cells.Item(81, 9).Value = 7599.75  # I81: 3459.4 -> 7599.75
$ws.Cells.Item(81, 11).Value = 15199.5  # K81: 6918.8 -> 15199.5
$ws.Cells.Item(81, 13).Value = -14138.5  # M81: -5857.8 -> -14138.5
$ws.Cells.Item(84, 8).Value = 6723  # H84: 4820.737 -> 6723
$ws.Cells.Item(84, 9).Value = 7599.75  # I84: 3459.4 -> 7599.75
$ws.Cells.Item(84, 11).Value = 75997.5  # K84: 34594 -> 75997.5
$ws.Cells.Item(84, 13).Value = -70693.5  # M84: -29290 -> -70693.5
$ws.Cells.Item(96, 8).Value = 7996.4287  # H96: 7996.5713 -> 7996.4287
$ws.Cells.Item(96, 9).Value = 7995.75  # I96: 7998 -> 7995.75
$ws.Cells.Item(96, 10).Value = 7997.3335  # J96: 7995.5 -> 7997.3335
$ws.Cells.Item(96, 11).Value = 7995.75  # K96: 7998 -> 7995.75
$ws.Cells.Item(96, 12).Value = 7997.3335  # L96: 7995.5 -> 7997.3335
$ws.Cells.Item(96, 13).Value = -6622.75  # M96: -6625 -> -6622.75
$ws.Cells.Item(96, 14).Value = -10743.3335  # N96: -10741.5 -> -10743.3335
$ws.Cells.Item(110, 8).Value = 60000  # H110: 0 -> 60000
$ws.Cells.Item(110, 10).Value = 60000  # J110: 0 -> 60000
$ws.Cells.Item(110, 12).Value = 60000  # L110: 0 -> 60000
$ws.Cells.Item(110, 14).Value = -68180  # N110: None -> -68180
$ws.Cells.Item(132, 8).Value = 3287.9868  # H132: 3187.9495 -> 3287.9868
$ws.Cells.Item(132, 9).Value = 3029.7014  # I132: 2955.8115 -> 3029.7014
$ws.Cells.Item(132, 10).Value = 5210.778  # J132: 4789.7 -> 5210.778
$ws.Cells.Item(132, 11).Value = 9089.1042  # K132: 8867.434499999999 -> 9089.1042
$ws.Cells.Item(132, 12).Value = 15632.334  # L132: 14369.1 -> 15632.334
$ws.Cells.Item(132, 13).Value = -6559.1042  # M132: -6337.434499999999 -> -6559.1042
$ws.Cells.Item(132, 14).Value = -20692.334  # N132: -19429.1 -> -20692.334

